$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 598
$ws.Range("I5").Value = 598
$ws.Range("K5").Value = 598
$ws.Range("M5").Value = -483
$ws.Range("H17").Value = 168159.83
$ws.Range("J17").Value = 183310.73
$ws.Range("L17").Value = 549932.1900000001
$ws.Range("N17").Value = -550268.1900000001
$ws.Range("H33").Value = 450.73334
$ws.Range("I33").Value = 450.73334
$ws.Range("K33").Value = 450.73334
$ws.Range("M33").Value = -221.73334
$ws.Range("H43").Value = 1645.875
$ws.Range("I43").Value = 1655.6
$ws.Range("J43").Value = 1629.6666
$ws.Range("K43").Value = 1655.6
$ws.Range("L43").Value = 1629.6666
$ws.Range("M43").Value = -1586.6
$ws.Range("N43").Value = -1767.6666
$ws.Range("H51").Value = 14199.4
$ws.Range("J51").Value = 16999.75
$ws.Range("L51").Value = 16999.75
$ws.Range("N51").Value = -17967.75
$ws.Range("H64").Value = 5000
$ws.Range("I64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("M64").Value = -4752
$ws.Range("H67").Value = 5000
$ws.Range("I67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("M67").Value = -4142
$ws.Range("H69").Value = 14288.842
$ws.Range("J69").Value = 16199.333
$ws.Range("L69").Value = 48597.999
$ws.Range("N69").Value = -50345.999
$ws.Range("H72").Value = 14288.842
$ws.Range("J72").Value = 16199.333
$ws.Range("L72").Value = 145793.997
$ws.Range("N72").Value = -154529.997
$ws.Range("H76").Value = 8573.583000000001
$ws.Range("I76").Value = 8197.5
$ws.Range("K76").Value = 8197.5
$ws.Range("M76").Value = -7882.5
$ws.Range("H79").Value = 8573.583000000001
$ws.Range("I79").Value = 8197.5
$ws.Range("K79").Value = 8197.5
$ws.Range("M79").Value = -7105.5
$ws.Range("H80").Value = 659.125
$ws.Range("I80").Value = 696.6
$ws.Range("J80").Value = 596.6667
$ws.Range("K80").Value = 2089.8
$ws.Range("L80").Value = 1790.0001
$ws.Range("M80").Value = -1091.8
$ws.Range("N80").Value = -3786.0001
$ws.Range("H83").Value = 659.125
$ws.Range("I83").Value = 696.6
$ws.Range("J83").Value = 596.6667
$ws.Range("K83").Value = 6269.400000000001
$ws.Range("L83").Value = 5370.0003
$ws.Range("M83").Value = -1277.400000000001
$ws.Range("N83").Value = -15354.0003
$ws.Range("H86").Value = 500000
$ws.Range("J86").Value = 500000
$ws.Range("L86").Value = 500000
$ws.Range("N86").Value = -502246
$ws.Range("H88").Value = 36916476
$ws.Range("J88").Value = 3981270
$ws.Range("L88").Value = 3981270
$ws.Range("N88").Value = -3982082
$ws.Range("H89").Value = 500000
$ws.Range("J89").Value = 500000
$ws.Range("L89").Value = 2500000
$ws.Range("N89").Value = -2511232
$ws.Range("H91").Value = 36916476
$ws.Range("J91").Value = 3981270
$ws.Range("L91").Value = 3981270
$ws.Range("N91").Value = -3984078
$ws.Range("H113").Value = 3125.75
$ws.Range("I113").Value = 2440
$ws.Range("J113").Value = 4268.6665
$ws.Range("K113").Value = 2440
$ws.Range("L113").Value = 4268.6665
$ws.Range("M113").Value = 814
$ws.Range("N113").Value = -10776.6665
$ws.Range("H115").Value = 339.66666
$ws.Range("J115").Value = 500
$ws.Range("L115").Value = 1500
$ws.Range("N115").Value = -4634
$ws.Range("H127").Value = 4212.125
$ws.Range("I127").Value = 3197
$ws.Range("J127").Value = 4357.143
$ws.Range("K127").Value = 9591
$ws.Range("L127").Value = 13071.429
$ws.Range("M127").Value = -4631
$ws.Range("N127").Value = -22991.429
$ws.Range("H129").Value = 2166.68
$ws.Range("I129").Value = 2899.5
$ws.Range("J129").Value = 2102.9565
$ws.Range("K129").Value = 8698.5
$ws.Range("L129").Value = 6308.869499999999
$ws.Range("M129").Value = -3698.5
$ws.Range("N129").Value = -16308.8695
$ws.Range("H138").Value = 3175.4468
$ws.Range("J138").Value = 3734.0688
$ws.Range("L138").Value = 11202.2064
$ws.Range("N138").Value = -21482.2064
$ws.Range("H141").Value = 2695.8333
$ws.Range("I141").Value = 2695.8333
$ws.Range("K141").Value = 8087.499899999999
$ws.Range("M141").Value = -2907.499899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 452761.9
$ws.Range("I2").Value = 926336.7
$ws.Range("K2").Value = 926336.7
$ws.Range("M2").Value = -926223.7
$ws.Range("H32").Value = 4527.5283
$ws.Range("I32").Value = 2454.238
$ws.Range("K32").Value = 2454.238
$ws.Range("M32").Value = -2167.238
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H45").Value = 1434.8889
$ws.Range("I45").Value = 1458
$ws.Range("J45").Value = 1250
$ws.Range("K45").Value = 1458
$ws.Range("L45").Value = 1250
$ws.Range("M45").Value = -1081
$ws.Range("N45").Value = -2004
$ws.Range("H61").Value = 33344998
$ws.Range("I61").Value = 34494276
$ws.Range("K61").Value = 34494276
$ws.Range("M61").Value = -34494064
$ws.Range("H74").Value = 58829770
$ws.Range("I74").Value = 71435500
$ws.Range("J74").Value = 2983.3333
$ws.Range("K74").Value = 71435500
$ws.Range("L74").Value = 2983.3333
$ws.Range("M74").Value = -71434626
$ws.Range("N74").Value = -4731.3333
$ws.Range("H77").Value = 58829770
$ws.Range("I77").Value = 71435500
$ws.Range("J77").Value = 2983.3333
$ws.Range("K77").Value = 357177500
$ws.Range("L77").Value = 14916.6665
$ws.Range("M77").Value = -357173132
$ws.Range("N77").Value = -23652.6665
$ws.Range("H110").Value = 129699.25
$ws.Range("I110").Value = 203916.6
$ws.Range("J110").Value = 6003.6665
$ws.Range("K110").Value = 203916.6
$ws.Range("L110").Value = 6003.6665
$ws.Range("M110").Value = -201871.6
$ws.Range("N110").Value = -10093.6665
$ws.Range("H116").Value = 452761.9
$ws.Range("I116").Value = 926336.7
$ws.Range("K116").Value = 926336.7
$ws.Range("M116").Value = -924042.7
$ws.Range("H122").Value = 1502.4333
$ws.Range("I122").Value = 1216.7037
$ws.Range("K122").Value = 3650.1111
$ws.Range("M122").Value = -1200.1111
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080
$ws.Range("H132").Value = 2044816.9
$ws.Range("I132").Value = 2944501
$ws.Range("J132").Value = 5532.6
$ws.Range("K132").Value = 8833503
$ws.Range("L132").Value = 16597.8
$ws.Range("M132").Value = -8830973
$ws.Range("N132").Value = -21657.8
$ws.Range("H136").Value = 33344998
$ws.Range("I136").Value = 34494276
$ws.Range("K136").Value = 103482828
$ws.Range("M136").Value = -103480278
$ws.Range("H140").Value = 59997.5
$ws.Range("J140").Value = 59997.5
$ws.Range("L140").Value = 59997.5
$ws.Range("N140").Value = -70357.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 452761.9
$ws.Range("I3").Value = 926336.7
$ws.Range("K3").Value = 926336.7
$ws.Range("M3").Value = -926222.7
$ws.Range("H105").Value = 4128.4287
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1815277.1
$ws.Range("I16").Value = 2176332.5
$ws.Range("K16").Value = 2176332.5
$ws.Range("M16").Value = -2176045.5
$ws.Range("H22").Value = 33666
$ws.Range("I22").Value = 99999
$ws.Range("J22").Value = 499.5
$ws.Range("K22").Value = 99999
$ws.Range("L22").Value = 499.5
$ws.Range("M22").Value = -99649
$ws.Range("N22").Value = -1199.5
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H31").Value = 12015.2
$ws.Range("I31").Value = 8426.608
$ws.Range("J31").Value = 16870.354
$ws.Range("K31").Value = 8426.608
$ws.Range("L31").Value = 16870.354
$ws.Range("M31").Value = -8131.608
$ws.Range("N31").Value = -17460.354
$ws.Range("H34").Value = 12015.2
$ws.Range("I34").Value = 8426.608
$ws.Range("J34").Value = 16870.354
$ws.Range("K34").Value = 8426.608
$ws.Range("L34").Value = 16870.354
$ws.Range("M34").Value = -8224.608
$ws.Range("N34").Value = -17274.354
$ws.Range("H86").Value = 4425.1113
$ws.Range("I86").Value = 4287
$ws.Range("K86").Value = 4287
$ws.Range("M86").Value = -3164
$ws.Range("H89").Value = 4425.1113
$ws.Range("I89").Value = 4287
$ws.Range("K89").Value = 21435
$ws.Range("M89").Value = -15819
$ws.Range("H92").Value = 29999.5
$ws.Range("J92").Value = 29999.5
$ws.Range("L92").Value = 29999.5
$ws.Range("N92").Value = -34991.5
$ws.Range("H107").Value = 541288.9
$ws.Range("I107").Value = 906336.4399999999
$ws.Range("K107").Value = 906336.4399999999
$ws.Range("M107").Value = -904416.4399999999
$ws.Range("H113").Value = 1815277.1
$ws.Range("I113").Value = 2176332.5
$ws.Range("K113").Value = 2176332.5
$ws.Range("M113").Value = -2174162.5
$ws.Range("H132").Value = 76925140
$ws.Range("I132").Value = 83335384
$ws.Range("J132").Value = 2199
$ws.Range("K132").Value = 250006152
$ws.Range("L132").Value = 6597
$ws.Range("M132").Value = -250003622
$ws.Range("N132").Value = -11657
$ws.Range("H134").Value = 8930484
$ws.Range("I134").Value = 10418339
$ws.Range("J134").Value = 3357.75
$ws.Range("K134").Value = 31255017
$ws.Range("L134").Value = 10073.25
$ws.Range("M134").Value = -31252482
$ws.Range("N134").Value = -15143.25
$ws.Range("H141").Value = 306242.75
$ws.Range("J141").Value = 344277.44
$ws.Range("L141").Value = 344277.44
$ws.Range("N141").Value = -354637.44

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 91.2
$ws.Range("I6").Value = 91.2
$ws.Range("K6").Value = 273.6
$ws.Range("M6").Value = -160.6
$ws.Range("H33").Value = 305.1905
$ws.Range("J33").Value = 433.91666
$ws.Range("L33").Value = 2603.49996
$ws.Range("N33").Value = -3169.49996
$ws.Range("H112").Value = 15960
$ws.Range("I112").Value = 9393.333000000001
$ws.Range("K112").Value = 28179.999
$ws.Range("M112").Value = -27071.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 391.55884
$ws.Range("I2").Value = 484.8
$ws.Range("K2").Value = 484.8
$ws.Range("M2").Value = -371.8
$ws.Range("H70").Value = 10661.429
$ws.Range("I70").Value = 10457.875
$ws.Range("J70").Value = 10932.833
$ws.Range("K70").Value = 10457.875
$ws.Range("L70").Value = 10932.833
$ws.Range("M70").Value = -10187.875
$ws.Range("N70").Value = -11472.833
$ws.Range("H73").Value = 10661.429
$ws.Range("I73").Value = 10457.875
$ws.Range("J73").Value = 10932.833
$ws.Range("K73").Value = 10457.875
$ws.Range("L73").Value = 10932.833
$ws.Range("M73").Value = -9521.875
$ws.Range("N73").Value = -12804.833
$ws.Range("H97").Value = 1235.9131
$ws.Range("I97").Value = 675.2857
$ws.Range("J97").Value = 2108
$ws.Range("K97").Value = 675.2857
$ws.Range("L97").Value = 2108
$ws.Range("M97").Value = -179.2857
$ws.Range("N97").Value = -3100
$ws.Range("H121").Value = 15000
$ws.Range("J121").Value = 15000
$ws.Range("L121").Value = 15000
$ws.Range("N121").Value = -18494
$ws.Range("H132").Value = 3128949
$ws.Range("I132").Value = 3574457.5
$ws.Range("J132").Value = 10390.8
$ws.Range("K132").Value = 10723372.5
$ws.Range("L132").Value = 31172.4
$ws.Range("M132").Value = -10720842.5
$ws.Range("N132").Value = -36232.39999999999
$ws.Range("H136").Value = 21361.191
$ws.Range("J136").Value = 21361.191
$ws.Range("L136").Value = 64083.573
$ws.Range("N136").Value = -69183.573

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2303.4614
$ws.Range("I16").Value = 699.4666999999999
$ws.Range("J16").Value = 4490.727
$ws.Range("K16").Value = 699.4666999999999
$ws.Range("L16").Value = 4490.727
$ws.Range("M16").Value = -529.4666999999999
$ws.Range("N16").Value = -4830.727
$ws.Range("H22").Value = 2810.2354
$ws.Range("I22").Value = 1379
$ws.Range("J22").Value = 4854.857
$ws.Range("K22").Value = 1379
$ws.Range("L22").Value = 4854.857
$ws.Range("M22").Value = -1084
$ws.Range("N22").Value = -5444.857
$ws.Range("H27").Value = 2810.2354
$ws.Range("I27").Value = 1379
$ws.Range("J27").Value = 4854.857
$ws.Range("K27").Value = 1379
$ws.Range("L27").Value = 4854.857
$ws.Range("M27").Value = -1272
$ws.Range("N27").Value = -5068.857
$ws.Range("H40").Value = 9666.666999999999
$ws.Range("I40").Value = 9000
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 9000
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = -8864
$ws.Range("N40").Value = -10272
$ws.Range("H46").Value = 1024.4615
$ws.Range("I46").Value = 742
$ws.Range("J46").Value = 1660
$ws.Range("K46").Value = 742
$ws.Range("L46").Value = 1660
$ws.Range("M46").Value = -554
$ws.Range("N46").Value = -2036
$ws.Range("H55").Value = 364.92856
$ws.Range("J55").Value = 473.625
$ws.Range("L55").Value = 473.625
$ws.Range("N55").Value = -819.625
$ws.Range("H68").Value = 2977550.5
$ws.Range("I68").Value = 2977550.5
$ws.Range("K68").Value = 2977550.5
$ws.Range("M68").Value = -2976801.5
$ws.Range("H71").Value = 2977550.5
$ws.Range("I71").Value = 2977550.5
$ws.Range("K71").Value = 14887752.5
$ws.Range("M71").Value = -14884008.5
$ws.Range("H82").Value = 1055.125
$ws.Range("I82").Value = 1090.3334
$ws.Range("K82").Value = 1090.3334
$ws.Range("M82").Value = -729.3334
$ws.Range("H85").Value = 1055.125
$ws.Range("I85").Value = 1090.3334
$ws.Range("K85").Value = 1090.3334
$ws.Range("M85").Value = 157.6666
$ws.Range("H95").Value = 34400
$ws.Range("J95").Value = 34400
$ws.Range("L95").Value = 34400
$ws.Range("N95").Value = -39892
$ws.Range("H104").Value = 2356.6667
$ws.Range("J104").Value = 2356.6667
$ws.Range("L104").Value = 2356.6667
$ws.Range("N104").Value = -9344.6667
$ws.Range("H132").Value = 36940316
$ws.Range("I132").Value = 36940316
$ws.Range("K132").Value = 110820948
$ws.Range("M132").Value = -110818418

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 16333
$ws.Range("I32").Value = 19999.5
$ws.Range("K32").Value = 19999.5
$ws.Range("M32").Value = -19682.5
$ws.Range("H41").Value = 15265
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 15265
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 15265
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -16045
$ws.Range("H62").Value = 25249.75
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 25249.75
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 25249.75
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -26497.75
$ws.Range("H65").Value = 25249.75
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 25249.75
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 126248.75
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -132488.75
$ws.Range("H81").Value = 3806.9285
$ws.Range("I81").Value = 4030.5
$ws.Range("J81").Value = 3248
$ws.Range("K81").Value = 8061
$ws.Range("L81").Value = 6496
$ws.Range("M81").Value = -7000
$ws.Range("N81").Value = -8618
$ws.Range("H84").Value = 3806.9285
$ws.Range("I84").Value = 4030.5
$ws.Range("J84").Value = 3248
$ws.Range("K84").Value = 40305
$ws.Range("L84").Value = 32480
$ws.Range("M84").Value = -35001
$ws.Range("N84").Value = -43088
$ws.Range("H122").Value = 6465.3335
$ws.Range("I122").Value = 7778.4443
$ws.Range("J122").Value = 2526
$ws.Range("K122").Value = 23335.3329
$ws.Range("L122").Value = 7578
$ws.Range("M122").Value = -20885.3329
$ws.Range("N122").Value = -12478
$ws.Range("H132").Value = 27788606
$ws.Range("I132").Value = 35720372
$ws.Range("J132").Value = 27424.5
$ws.Range("K132").Value = 107161116
$ws.Range("L132").Value = 82273.5
$ws.Range("M132").Value = -107158586
$ws.Range("N132").Value = -87333.5
$ws.Range("H136").Value = 100002060
$ws.Range("I136").Value = 125001190
$ws.Range("K136").Value = 375003570
$ws.Range("M136").Value = -375001020
